$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; temporarily unprotect to write values, then restore protection.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure text (A16).
$ws.Cells.Replace("2021-07-07", "2021-07-08") | Out-Null
# Replace() can trigger a row autofit side effect on the multi-line cell; re-autofit
# to keep the row height identical to its original (default) state.
$ws.Rows(16).AutoFit()

# Updated Weight (D) / Percent Change (E) values for rows 2-13.
$newValues = @{
    2  = @(0.02705781117107306, -0.002907915993538124)
    3  = @(0.02143313808958444, -0.003814064362336222)
    4  = @(0.05650991192515847, -0.002350728725905138)
    5  = @(0.1395722773329422, -0.003091190108191699)
    6  = @(0.02001993046250958, -0.002812939521800284)
    7  = @(0.1286206278356828, -0.004019978072846819)
    8  = @(0.08910463129495892, -0.01357798165137625)
    9  = @(0.02939998680688308, -0.01513840830449842)
    10 = @(0.1025316365153232, -0.0200385356454722)
    11 = @(0.2963674762191797, -0.009892759165350484)
    12 = @(0.08938257234670467, -0.01042619352478502)
    13 = @(1, -0.008871406053597042)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

$ws.Protect("D382")
